$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column C for "User_fb_url" (shifts old Comment column C -> D) ---
$ws.Columns.Item(3).Insert()

# --- Header row ---
$ws.Range("C1").Value = "User_fb_url"

# Header C1 ("User_fb_url") gets the same bold+border+center/top style as B1 used to have
# and as D1 ("Comment") still has.
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4160

# Header B1 ("User Name") keeps bold+border+center but loses the vertical=top alignment.
$ws.Range("B1").VerticalAlignment = -4107

# --- New column C data (Facebook profile urls); rows 7 and 10 stay empty, mirroring
#     the gaps that already existed in the Comment column for those rows. ---
$ws.Range("C2").Value = "/ammar.laabidi.52?rc=p&__tn__=R"
$ws.Range("C3").Value = "/mstirsalahn1?rc=p&__tn__=R"
$ws.Range("C4").Value = "/mouhammedamin.touati?rc=p&__tn__=R"
$ws.Range("C5").Value = "/wejdane.jedaydie?rc=p&__tn__=R"
$ws.Range("C6").Value = "/profile.php?id=100057178865281&rc=p&__tn__=R"
$ws.Range("C8").Value = "/nizar.hamdi.9?rc=p&__tn__=R"
$ws.Range("C9").Value = "/sinen.mh.33?rc=p&__tn__=R"

# --- Column B (User Name) data cells get centered alignment (new style) ---
$ws.Range("B2:B10").HorizontalAlignment = -4108

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 23.67
$ws.Columns.Item(3).ColumnWidth = 44.5
$ws.Columns.Item(4).ColumnWidth = 54.33

# --- Active selection as in the saved file ---
$ws.Range("C8").Select() | Out-Null
